## Apply the "Distance Ratio Statistics" update:
##  - Add the row-8 (minPts=8) statistics to the "All" and "Without Outliers" sheets
##  - A few previously-#VALUE!/#NUM! error cells on rows 3-7 now evaluate to #N/A
##  - The "Ratios" sheet recomputes automatically from the 'Without Outliers'/All formulas
##  - Update the remembered cell selection (row 7 -> row 8) on all three sheets

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "All"
# ---------------------------------------------------------------------------
$all = $wb.Worksheets.Item("All")

# Rows 3-5: Mode (D) / Kurtosis (J) now resolve to #N/A instead of #VALUE!/#NUM!
$all.Range("D3").Value = "#N/A"
$all.Range("J3").Value = "#N/A"
$all.Range("D4").Value = "#N/A"
$all.Range("J4").Value = "#N/A"
$all.Range("D5").Value = "#N/A"

# Row 8 (minPts = 8) statistics
$all.Range("B8").Value = 28.8199233333333
$all.Range("C8").Value = 9.58607666666667
$all.Range("D8").Value = 38.406
$all.Range("E8").Value = 38.406
$all.Range("F8").Value = 24.026885
$all.Range("G8").Value = 38.406
$all.Range("H8").Value = 275.678597577633
$all.Range("I8").Value = 16.6035718319172
$all.Range("J8").Value = "#NUM!"
$all.Range("K8").Value = -1.73205080756888
$all.Range("L8").Value = 28.75823
$all.Range("M8").Value = 9.64777
$all.Range("N8").Value = 38.406
$all.Range("O8").Value = 86.45977
$all.Range("P8").Value = 3

# ---------------------------------------------------------------------------
# Sheet "Without Outliers"
# ---------------------------------------------------------------------------
$wo = $wb.Worksheets.Item("Without Outliers")

# Rows 3-7: Mode (D) / Kurtosis (J) / Skewness (K) now resolve to #N/A
$wo.Range("D3").Value = "#N/A"
$wo.Range("J3").Value = "#N/A"
$wo.Range("K3").Value = "#N/A"
$wo.Range("D4").Value = "#N/A"
$wo.Range("J4").Value = "#N/A"
$wo.Range("K4").Value = "#N/A"
$wo.Range("D5").Value = "#N/A"
$wo.Range("D6").Value = "#N/A"
$wo.Range("J6").Value = "#N/A"
$wo.Range("D7").Value = "#N/A"
$wo.Range("J7").Value = "#N/A"

# Row 8 (minPts = 8) statistics
$wo.Range("B8").Value = 24.026885
$wo.Range("C8").Value = 14.379115
$wo.Range("D8").Value = "#VALUE!"
$wo.Range("E8").Value = 24.026885
$wo.Range("F8").Value = 16.8373275
$wo.Range("G8").Value = 31.2164425
$wo.Range("H8").Value = 413.51789636645
$wo.Range("I8").Value = 20.3351394479224
$wo.Range("J8").Value = "#NUM!"
$wo.Range("K8").Value = "#NUM!"
$wo.Range("L8").Value = 28.75823
$wo.Range("M8").Value = 9.64777
$wo.Range("N8").Value = 38.406
$wo.Range("O8").Value = 48.05377
$wo.Range("P8").Value = 2

# ---------------------------------------------------------------------------
# Sheet "Ratios" -- values are derived via ='Without Outliers'!x/All!x formulas,
# so they recalculate automatically now that row 8 has real data on both
# source sheets.
# ---------------------------------------------------------------------------
$ratios = $wb.Worksheets.Item("Ratios")
$excel.Calculate()

# ---------------------------------------------------------------------------
# Update the remembered cell selection (row 7 -> row 8) on all three sheets.
# "Without Outliers" is selected last so it remains the workbook's active
# (tabSelected) sheet, matching the original file.
# ---------------------------------------------------------------------------
$all.Range("B8:P8").Select() | Out-Null
$ratios.Range("B8:P8").Select() | Out-Null
$wo.Range("B8:P8").Select() | Out-Null
